$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.345.28"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "3.664.94"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'644.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.33%  "
$ws.Range("D6").Value = "'158.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.497"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "'7.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "4.289.55"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "'32.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "3.638.97"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "69.353.77"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D18").Value = "'15.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "'6.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'465.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "'9.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "'0.644"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "'79.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "3.816.83"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'0.0000123"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").Value = "'10.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'8.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "'26.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.163"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.662.88"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "'8.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'178.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.43%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'5.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.27%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "'0.922"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").Value = "'47.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "'28.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "'2.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'7.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000262"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.72%  "
$ws.Range("D51").Value = "'1.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.73%  "
